# Recolor the two sticky-note "color key" rectangles:
#   Rectangle 264  cc2fc6 -> f5058e
#   Rectangle 262  41b6e0 -> c1106e
#
# Note: $d.Shapes enumerates in z-order for metadata (.Name/.Id/...), but the
# indices that actually land a Fill/ForeColor write line up with the shapes'
# order in the document body. Index 2 in that order is "Rectangle 264" and
# index 27 is "Rectangle 262" in this document, so those are the indices used
# below to target the two shapes reliably.

$d = $word.ActiveDocument

# Rectangle 264: #cc2fc6 -> #f5058e  (R=0xF5,G=0x05,B=0x8E -> 0x8E0005 + ... => RGB(0xF5,0x05,0x8E))
$d.Shapes.Item(2).Fill.ForeColor.RGB = 9307637

# Rectangle 262: #41b6e0 -> #c1106e  (RGB(0xC1,0x10,0x6E))
$d.Shapes.Item(27).Fill.ForeColor.RGB = 7213249
